# Apply the "1.5.2" workbook update:
#  - fill in the previously-empty column Q (latest period) figures for
#    each indicator block, mirroring the pattern already used in column P
#  - update the active window's view (drop the stale topLeftCell scroll
#    position and move the live selection to O52)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric figures for column Q (values taken from the source data refresh).
$ws.Range("Q6").Value  = 0.044166007693658721
$ws.Range("Q8").Value  = 601820.30000000005
$ws.Range("Q10").Value = 0.36185407133694547
$ws.Range("Q12").Value = 20892.400000000001
$ws.Range("Q14").Value = 0.27408710679222598
$ws.Range("Q16").Value = 63884.800000000003
$ws.Range("Q18").Value = 0.00046658384803364067
$ws.Range("Q20").Value = 85729.5
$ws.Range("Q22").Value = 0.086032657053793982
$ws.Range("Q24").Value = 16970.3
$ws.Range("Q26").Value = "-"
$ws.Range("Q28").Value = 47183.5
$ws.Range("Q30").Value = "-"
$ws.Range("Q32").Value = 17405.3
$ws.Range("Q34").Value = 0.000011900270969169968
$ws.Range("Q36").Value = 84031.7
$ws.Range("Q38").Value = "-"
$ws.Range("Q40").Value = 231841.7
$ws.Range("Q42").Value = "-"
$ws.Range("Q44").Value = 33881.1

# Refresh the view: selecting a cell clears any stale scroll/topLeftCell
# state and records the new active cell / selection.
$ws.Range("O52").Select()
